$d = $word.ActiveDocument

# 1. Add a first-line indent (0.5in = 36pt = 720 twips) to the title paragraph
#    "Persuasive Paper Reading Journal"
$p1 = $d.Paragraphs(1)
$p1.Format.FirstLineIndent = 36

# 2. Merge the three runs that make up the "To read your selected novel ..."
#    sentence (removing the gramStart/gramEnd proofErr wrapping "in order to")
#    by replacing the full sentence text with itself.
$old1 = "To read your selected novel in order to write a persuasive paper which answers the question- what message is the author attempting to convey to his/her readers and how?  "
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $old1, 2) | Out-Null

# 3. Merge the three runs that make up the " - As you read, highlight and label..."
#    sentence (removing the gramStart/gramEnd proofErr wrapping "highlight")
#    by replacing the full sentence text with itself.
$old2 = " – As you read, highlight and label significant sentences from the story that reveal something interesting or important about a character, setting, conflict, or other notable items like the author’s use of figurative language.  This can be done easily in "
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $old2, 2) | Out-Null
